$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage (matches original inlineStr type) for D-column cells
# whose new value would otherwise be auto-parsed by Excel as a number.
$textCells = @('D4', 'D5', 'D6', 'D13', 'D14', 'D20', 'D21', 'D24', 'D25', 'D26', 'D27', 'D30', 'D33', 'D37', 'D40', 'D41', 'D42', 'D43', 'D44', 'D45', 'D48')
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated cell values exactly as captured in the diff.
$ws.Range('D2').Value = '61.662.06'
$ws.Range('E2').Value = '  -1.79%  '
$ws.Range('D3').Value = '3.001.61'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '597.85'
$ws.Range('E5').Value = '  +2.43%  '
$ws.Range('D6').Value = '144.59'
$ws.Range('E6').Value = '  -3.41%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('E8').Value = '  -0.69%  '
$ws.Range('D9').Value = '3.001.36'
$ws.Range('E9').Value = '  -0.97%  '
$ws.Range('E10').Value = '  -2.70%  '
$ws.Range('E11').Value = '  +4.12%  '
$ws.Range('E12').Value = '  +4.39%  '
$ws.Range('D13').Value = '0.0000229'
$ws.Range('E13').Value = '  -1.21%  '
$ws.Range('D14').Value = '34.37'
$ws.Range('E14').Value = '  -3.09%  '
$ws.Range('E15').Value = '  +2.70%  '
$ws.Range('D16').Value = '3.495.43'
$ws.Range('E16').Value = '  -1.13%  '
$ws.Range('E17').Value = '  -0.82%  '
$ws.Range('D18').Value = '61.581.07'
$ws.Range('D19').Value = '2.997.87'
$ws.Range('E19').Value = '  -1.16%  '
$ws.Range('D20').Value = '455.66'
$ws.Range('E20').Value = '  -2.72%  '
$ws.Range('D21').Value = '14.06'
$ws.Range('E21').Value = '  -0.17%  '
$ws.Range('E22').Value = '  -0.50%  '
$ws.Range('E23').Value = '  -0.73%  '
$ws.Range('D24').Value = '82.41'
$ws.Range('E24').Value = '  +1.56%  '
$ws.Range('D25').Value = '2.21'
$ws.Range('E25').Value = '  -8.10%  '
$ws.Range('D26').Value = '12.22'
$ws.Range('E26').Value = '  -1.69%  '
$ws.Range('D27').Value = '10.48'
$ws.Range('E27').Value = '  -0.41%  '
$ws.Range('E28').Value = '  +0.02%  '
$ws.Range('E29').Value = '  +1.77%  '
$ws.Range('D30').Value = '1.00'
$ws.Range('E30').Value = '  -0.02%  '
$ws.Range('E31').Value = '  -3.23%  '
$ws.Range('E32').Value = '  -4.22%  '
$ws.Range('D33').Value = '27.29'
$ws.Range('E33').Value = '  -0.72%  '
$ws.Range('E34').Value = '  -0.48%  '
$ws.Range('D35').Value = '0.0₃0821'
$ws.Range('E35').Value = '  +2.73%  '
$ws.Range('E36').Value = '  -1.64%  '
$ws.Range('D37').Value = '5.77'
$ws.Range('E37').Value = '  -0.37%  '
$ws.Range('E38').Value = '  -3.33%  '
$ws.Range('E39').Value = '  +1.86%  '
$ws.Range('D40').Value = '50.33'
$ws.Range('E40').Value = '  +0.07%  '
$ws.Range('B41').Value = 'dogwifhat'
$ws.Range('C41').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D41').Value = '2.89'
$ws.Range('E41').Value = '  -2.97%  '
$ws.Range('B42').Value = 'Kaspa'
$ws.Range('C42').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D42').Value = '0.122'
$ws.Range('E42').Value = '  +8.60%  '
$ws.Range('D43').Value = '400.80'
$ws.Range('E43').Value = '  -5.45%  '
$ws.Range('D44').Value = '39.58'
$ws.Range('E44').Value = '  +4.50%  '
$ws.Range('D45').Value = '0.0354'
$ws.Range('E45').Value = '  -0.84%  '
$ws.Range('E46').Value = '  -5.49%  '
$ws.Range('D47').Value = '2.721.58'
$ws.Range('E47').Value = '  -2.90%  '
$ws.Range('D48').Value = '133.17'
$ws.Range('E48').Value = '  +2.35%  '
$ws.Range('E49').Value = '  +0.08%  '
$ws.Range('E50').Value = '  +1.54%  '
$ws.Range('E51').Value = '  -0.56%  '
